$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text formatting so that
# numeric-looking values (e.g. "97.50") are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.316.81"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.269.75"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "305.69"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "97.50"
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "35.73"
$ws.Range("E10").Value = "  +8.94%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "6.65"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "2.599.98"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "14.43"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "2.268.81"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "42.222.42"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "12.57"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "67.71"
$ws.Range("D23").Value = "237.76"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "37.43"
$ws.Range("E28").Value = "  +6.28%  "
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").Value = "160.48"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "3.17"
$ws.Range("E34").Value = "  +4.98%  "
$ws.Range("D35").Value = "0.0742"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").Value = "17.22"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "4.07"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").Value = "2.45"
$ws.Range("E42").Value = "  +14.61%  "
$ws.Range("D43").Value = "1.993.08"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "19.01"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "9.97"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").Value = "53.55"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "1.52"
$ws.Range("D50").Value = "72.03"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "91.58"
$ws.Range("E51").Value = "  -0.28%  "
